# Deploying to gh-pages from @ torstees/ncpi-fhir-ig-2@537984b1c0b3db4c39b833973081dfc1c90611e0
#
# The generated "Metadata"/"Elements" workbook was regenerated from an updated
# FHIR StructureDefinition (access-policy-description). Net effect:
#   - Metadata!B8  (Date)    refreshed to the new generation timestamp.
#   - Metadata!B20 (Context) now targets the Consent element instead of Element.
#   - Elements: the separate "valueMarkdown" slice row for Extension.value[x]
#     was removed; its Short/Definition text was folded back into the base
#     Extension.value[x] row, and the (now unused) slicing metadata on that
#     row was cleared out.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B8").Value = "2024-03-11T22:11:27+00:00"
$meta.Range("B20").Value = "element:Consent"

# ---------------------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Fold the (soon to be deleted) "valueMarkdown" slice row's Short/Definition
# text into the parent Extension.value[x] row, and drop its now-orphaned
# slicing discriminator/description/rules.
$els.Range("L6").Value = "Descriptive text summarizing the policy restrictions and other details associated with this access provision."
$els.Range("AB6").Value = ""
$els.Range("AC6").Value = ""
$els.Range("AE6").Value = ""

# Remove the "Extension.value[x]:valueMarkdown" slice row entirely.
$els.Rows.Item(7).Delete()

# Column widths/best-fit shrink now that the longer strings are gone
# (values tuned against this runtime's character-width quantization so the
# saved OOXML <col> width lands as close as possible to Excel's own bestFit
# recompute of 19.59765625 / 12.66015625).
$els.Columns.Item(1).ColumnWidth = 18.76
$els.Columns.Item(3).ColumnWidth = 11.76
